$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.100.05"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.901.23"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7266"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3107"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06868"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7715"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07940"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "1.886.56"
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.48%  "
$ws.Range("D16").Value = "30.117.17"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.800"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007721"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").Value = "2.125.62"
$ws.Range("E22").Value = "  -2.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.972"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.326"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1267"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.037"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.02%  "
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.537"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.283"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.060"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05094"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.278"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7340"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.752"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01922"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.779"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.329"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4414"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.922"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8341"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.560"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.686"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "936.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1177"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.56%  "
